# Auto-generated edit script: update cryptos price/volume data per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.893.64'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.630.37'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.519'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.44'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0881'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '1.861.84'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '1.628.16'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('E14').Value = '  -1.34%  '
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('D17').Value = '27.896.33'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.29'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.88'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').Value = '1.389.34'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  +10.59%  '
$ws.Range('E37').Value = '  -0.71%  '
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.557'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.848'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.60'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').Value = '1.772.22'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  -2.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.49'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.64'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.43%  '
